# Apply updated crypto price/volume data to the worksheet.
# Source: scheduled GitHub Actions scrape update (cryptos list).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new Price text, new Volume(1h) text).
# Price values that are purely numeric-looking (e.g. "1.017") are written with a
# leading apostrophe so Excel stores them as literal text (matching the source
# data, which keeps things like trailing zeros, e.g. "5.950") instead of coercing
# them to numbers.

$ws.Range("D2").Value = '27.903.17'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '1.878.77'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("D4").Value = '''1.017'
$ws.Range("E4").Value = '  +1.40%  '
$ws.Range("D5").Value = '''334.95'
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").Value = '''1.016'
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("E7").Value = '  -0.26%  '
$ws.Range("D8").Value = '''0.3911'
$ws.Range("E8").Value = '  -1.30%  '
$ws.Range("D9").Value = '''46.82'
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("D10").Value = '''0.07955'
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").Value = '''1.008'
$ws.Range("E11").Value = '  -1.49%  '
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = '1.870.14'
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("D14").Value = '''5.950'
$ws.Range("E14").Value = '  -0.27%  '
$ws.Range("D15").Value = '''7.101'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("E16").Value = '  +1.54%  '
$ws.Range("D17").Value = '''0.06785'
$ws.Range("E17").Value = '  +2.62%  '
$ws.Range("D18").Value = '''87.49'
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("E20").Value = '  -1.13%  '
$ws.Range("D21").Value = '''1.016'
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D22").Value = '27.897.53'
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("D23").Value = '''5.474'
$ws.Range("E23").Value = '  -0.63%  '
$ws.Range("D24").Value = '''10.96'
$ws.Range("E24").Value = '  -0.59%  '
$ws.Range("D25").Value = '''2.358'
$ws.Range("E25").Value = '  +2.50%  '
$ws.Range("D26").Value = '2.113.17'
$ws.Range("E26").Value = '  +0.84%  '
$ws.Range("D27").Value = '''159.79'
$ws.Range("E27").Value = '  +2.02%  '
$ws.Range("D28").Value = '''19.91'
$ws.Range("E28").Value = '  -1.77%  '
$ws.Range("D29").Value = '''2.081'
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("D30").Value = '''5.449'
$ws.Range("E30").Value = '  -2.42%  '
$ws.Range("D31").Value = '''120.86'
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("D32").Value = '''0.09536'
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("D33").Value = '''0.9605'
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("D34").Value = '''3.653'
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("D35").Value = '''5.316'
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").Value = '''1.348'
$ws.Range("E36").Value = '  -7.24%  '
$ws.Range("D37").Value = '''0.06110'
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").Value = '''0.02242'
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("E39").Value = '  -1.98%  '
$ws.Range("D40").Value = '''1.014'
$ws.Range("E40").Value = '  +1.19%  '
$ws.Range("D41").Value = '''8.113'
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("D42").Value = '''0.5908'
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").Value = '''0.1892'
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("D44").Value = '''10.21'
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("D45").Value = '''1.269'
$ws.Range("E45").Value = '  +0.62%  '
$ws.Range("D46").Value = '''0.5649'
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("D47").Value = '''12.22'
$ws.Range("E47").Value = '  -0.47%  '
$ws.Range("D48").Value = '''3.393'
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("D49").Value = '''1.919'
$ws.Range("E49").Value = '  -0.72%  '
$ws.Range("D50").Value = '''0.06859'
$ws.Range("E50").Value = '  +0.53%  '
$ws.Range("D51").Value = '''113.50'
$ws.Range("E51").Value = '  +1.30%  '
